# Daily attendance processing - 2026-02-07 08:00:05 UTC
# Swap "Miss Dina Nasr, Administrator" -> "Administrator, Miss Dina Nasr"
# in column G of the active worksheet (every row where the exact text matches).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "Miss Dina Nasr, Administrator"
$newText = "Administrator, Miss Dina Nasr"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
